$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "See San Antonio three ways with this combination bus, boat, and viewing tower experience. "
$ws.Range("C3").Value = "Visit five missions filled with history and religious importance."
$ws.Range("C4").Value = "Tour some of the most haunted sites in the city at night. Guide leads the way by lantern as you hear about the paranormal activity. "
$ws.Range("C5").Value = "Home to thousands of exotic animals, guaranteed hit with people of all ages. "
$ws.Range("C6").Value = "Cruise through San Antonio on a 2-hour Segway tour that takes in the highlights of historic downtown."
